$wb = $excel.ActiveWorkbook

$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

$zh.Range("E4").Value = "2016-03-21 04:41:37"
$zh.Range("H4").Value = "2016-03-21 04:41:59"

$de.Range("E4").Value = "2016-03-21 04:41:40"
$de.Range("H4").Value = "2016-03-21 04:42:05"
